$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Protokoll")

# New comment text for 3.Sprint (C7) referencing the new shared string
$ws.Range("C7").Value = "Anforderung von Key für die API, erste MockUps und Grunddesign AndroidApp"

# Progress values for 3.Sprint (row 8): 100% done tasks, 50% overall progress
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0.5

# Restore the active cell selection as recorded after the edit
$ws.Range("D12").Select() | Out-Null
